$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1) Make a throw-away duplicate of "ODI Batting" purely so we have a sheet
#    whose cells are "clipboard eligible" for a cross-sheet style copy (a
#    sheet as freshly loaded from disk can't be copy/pasted from directly in
#    this host). This lets the new sheet's header reuse the exact same bold/
#    border/centered style entry as ODI Batting's header instead of a newly
#    synthesized (merely similar) style.
# ---------------------------------------------------------------------------
$odiBattingOrig = $wb.Worksheets.Item("ODI Batting")
$odiBattingOrig.Copy($odiBattingOrig, $null)
$scratch = $wb.Worksheets.Item(1)
$scratch.Name = "Scratch"

# .Copy() rebinds the variable used to call it to the *new* sheet, so fetch a
# fresh handle to the real "ODI Batting" sheet by name.
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------------
# 2) Insert the real new "Player Info" worksheet right before "ODI Batting"
#    and stamp its header with the copied formatting.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$scratch.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

# Discard the scratch helper sheet now - it was never part of the intended
# result. (Deleting a sheet can invalidate other live worksheet references in
# this host, so we re-fetch everything we still need by name afterwards.)
$scratch.Delete() | Out-Null

$playerInfo = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Header captions
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row
$playerInfo.Range("A2").Value = "'5938"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Kjorn Yohance Ottley"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ---------------------------------------------------------------------------
# 3) On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and collapse the
#    full scorecard URLs down to just the bare numeric match code.
# ---------------------------------------------------------------------------
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Value = "'4445"
$odiBatting.Range("D2").Style = "Normal"
$odiBatting.Range("D3").Value = "'4447"
$odiBatting.Range("D3").Style = "Normal"

# Match the original workbook's "first sheet is active" view state.
$playerInfo.Activate()
